$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 2, shifting everything else (and row 1 header stays) down
$ws.Rows.Item(2).Insert()

# Row 3 (the old row 2) already carries the "quotePrefix" style used for
# wave numbers like "96.3" -- copy that formatting onto the new B2 so we
# reuse the existing style entry instead of minting a new one.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill the new row with the new Eurobarometer wave entry (order chosen to
# match shared-string append order: archive_id, wave, description, timeframe)
$ws.Range("A2").Value = "ZA7886"
$ws.Range("B2").Value = "'97.1"
$ws.Range("D2").Value = "Europeans, Agriculture and the CAP (COVID-19 Pandemic)"
$ws.Range("C2").Value = "February-March 2022"

# Update the view/selection to match
$ws.Range("C2").Select()
